$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the Street value in rows 4 and 5: "Redwood Drive" -> "2 Redwood Drive"
$ws.Range("F4").Value = "2 Redwood Drive"
$ws.Range("F5").Value = "2 Redwood Drive"

# Remove the stray test rows (6-9) that held leftover/no-longer-needed data
$ws.Range("A6:AB9").EntireRow.Delete()

# Update the active selection to match the post-cleanup state
[void]$ws.Range("A6:XFD63").Select()
